# Add the "2022-Q3" quarterly holdings sheet and record it in the "总计"
# (totals) summary sheet.
#
# Source workbook layout (before):
#   总计, 2022-Q2, 2022-Q1, 2021-Q4, 2021-Q3, 2021-Q2, 2021-Q1, 2020-Q4
#
# Target layout (after):
#   总计, 2022-Q3, 2022-Q2, 2022-Q1, 2021-Q4, 2021-Q3, 2021-Q2, 2021-Q1, 2020-Q4

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert the new "2022-Q3" worksheet right after "总计" (i.e. before the
#    sheet that is currently in position 2, "2022-Q2").
# ---------------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item(1)
$q2Sheet = $wb.Worksheets.Item(2)
$ws = $wb.Worksheets.Add($q2Sheet)
$ws.Name = "2022-Q3"

# Style template cells (already formatted with the bold/centered/bordered
# "index" style used throughout this workbook) that we reuse via
# Copy + PasteSpecial(xlPasteFormats) so no new style entries are created.
$headerStyleSrc = $q2Sheet.Range("B1")
$indexStyleSrc = $q2Sheet.Range("A2")

# Headers
$ws.Range("B1").Value = "基金代码"
$ws.Range("C1").Value = "基金名称"
$ws.Range("D1").Value = "基金规模"
$ws.Range("E1").Value = "股票总仓位"
$ws.Range("F1").Value = "仓位占比"
$ws.Range("G1").Value = "持有市值(亿元)"
$ws.Range("H1").Value = "仓位排名"

$headerStyleSrc.Copy()
$ws.Range("B1:H1").PasteSpecial(-4122)

# Data rows (2022-Q3 fund holdings)
$ws.Cells.Item(2,1).Value = 0
$ws.Cells.Item(2,2).Value = "003745"
$ws.Cells.Item(2,3).Value = "广发多元新兴股票"
$ws.Cells.Item(2,4).Value = "35.14"
$ws.Cells.Item(2,5).Value = "90.17"
$ws.Cells.Item(2,6).Value = "6.69"
$ws.Cells.Item(2,7).Value = "2.3509"
$ws.Cells.Item(2,8).Value = 6

$ws.Cells.Item(3,1).Value = 1
$ws.Cells.Item(3,2).Value = "001239"
$ws.Cells.Item(3,3).Value = "长盛国企改革主题灵活配置混合"
$ws.Cells.Item(3,4).Value = "4.46"
$ws.Cells.Item(3,5).Value = "90.97"
$ws.Cells.Item(3,6).Value = "9.04"
$ws.Cells.Item(3,7).Value = "0.4032"
$ws.Cells.Item(3,8).Value = 1

$ws.Cells.Item(4,1).Value = 2
$ws.Cells.Item(4,2).Value = "501201"
$ws.Cells.Item(4,3).Value = "红土创新科技创新 3 年封闭运作灵活配置混合"
$ws.Cells.Item(4,4).Value = "3.85"
$ws.Cells.Item(4,5).Value = "98.34"
$ws.Cells.Item(4,6).Value = "3.87"
$ws.Cells.Item(4,7).Value = "0.1490"
$ws.Cells.Item(4,8).Value = 7

$ws.Cells.Item(5,1).Value = 3
$ws.Cells.Item(5,2).Value = "006449"
$ws.Cells.Item(5,3).Value = "浙商汇金量化精选灵活配置混合"
$ws.Cells.Item(5,4).Value = "2.57"
$ws.Cells.Item(5,5).Value = "81.91"
$ws.Cells.Item(5,6).Value = "3.06"
$ws.Cells.Item(5,7).Value = "0.0786"
$ws.Cells.Item(5,8).Value = 10

$ws.Cells.Item(6,1).Value = 4
$ws.Cells.Item(6,2).Value = "015071"
$ws.Cells.Item(6,3).Value = "鑫元专精特新混合A"
$ws.Cells.Item(6,4).Value = "2.65"
$ws.Cells.Item(6,5).Value = "74.01"
$ws.Cells.Item(6,6).Value = "2.22"
$ws.Cells.Item(6,7).Value = "0.0588"
$ws.Cells.Item(6,8).Value = 7

$ws.Cells.Item(7,1).Value = 5
$ws.Cells.Item(7,2).Value = "009432"
$ws.Cells.Item(7,3).Value = "德邦科技创新一年定期开放混合A"
$ws.Cells.Item(7,4).Value = "1.47"
$ws.Cells.Item(7,5).Value = "85.46"
$ws.Cells.Item(7,6).Value = "3.96"
$ws.Cells.Item(7,7).Value = "0.0582"
$ws.Cells.Item(7,8).Value = 5

$ws.Cells.Item(8,1).Value = 6
$ws.Cells.Item(8,2).Value = "970043"
$ws.Cells.Item(8,3).Value = "东吴裕盈一年持有期灵活配置混合A"
$ws.Cells.Item(8,4).Value = "0.96"
$ws.Cells.Item(8,5).Value = "52.43"
$ws.Cells.Item(8,6).Value = "5.52"
$ws.Cells.Item(8,7).Value = "0.0530"
$ws.Cells.Item(8,8).Value = 3

$ws.Cells.Item(9,1).Value = 7
$ws.Cells.Item(9,2).Value = "168401"
$ws.Cells.Item(9,3).Value = "红土创新转型精选灵活配置混合（LOF）"
$ws.Cells.Item(9,4).Value = "0.87"
$ws.Cells.Item(9,5).Value = "92.08"
$ws.Cells.Item(9,6).Value = "3.88"
$ws.Cells.Item(9,7).Value = "0.0338"
$ws.Cells.Item(9,8).Value = 6

$ws.Cells.Item(10,1).Value = 8
$ws.Cells.Item(10,2).Value = "000354"
$ws.Cells.Item(10,3).Value = "长盛城镇化主题混合"
$ws.Cells.Item(10,4).Value = "0.34"
$ws.Cells.Item(10,5).Value = "92.91"
$ws.Cells.Item(10,6).Value = "7.81"
$ws.Cells.Item(10,7).Value = "0.0266"
$ws.Cells.Item(10,8).Value = 2

$ws.Cells.Item(11,1).Value = 9
$ws.Cells.Item(11,2).Value = "970045"
$ws.Cells.Item(11,3).Value = "东吴裕盈一年持有期灵活配置混合C"
$ws.Cells.Item(11,4).Value = "0.44"
$ws.Cells.Item(11,5).Value = "52.43"
$ws.Cells.Item(11,6).Value = "5.52"
$ws.Cells.Item(11,7).Value = "0.0243"
$ws.Cells.Item(11,8).Value = 3

$ws.Cells.Item(12,1).Value = 10
$ws.Cells.Item(12,2).Value = "970044"
$ws.Cells.Item(12,3).Value = "东吴裕盈一年持有期灵活配置混合B"
$ws.Cells.Item(12,4).Value = "0.27"
$ws.Cells.Item(12,5).Value = "52.43"
$ws.Cells.Item(12,6).Value = "5.52"
$ws.Cells.Item(12,7).Value = "0.0149"
$ws.Cells.Item(12,8).Value = 3

$ws.Cells.Item(13,1).Value = 11
$ws.Cells.Item(13,2).Value = "009433"
$ws.Cells.Item(13,3).Value = "德邦科技创新一年定期开放混合C"
$ws.Cells.Item(13,4).Value = "0.36"
$ws.Cells.Item(13,5).Value = "85.46"
$ws.Cells.Item(13,6).Value = "3.96"
$ws.Cells.Item(13,7).Value = "0.0143"
$ws.Cells.Item(13,8).Value = 5

$ws.Cells.Item(14,1).Value = 12
$ws.Cells.Item(14,2).Value = "015072"
$ws.Cells.Item(14,3).Value = "鑫元专精特新混合C"
$ws.Cells.Item(14,4).Value = "0.25"
$ws.Cells.Item(14,5).Value = "74.01"
$ws.Cells.Item(14,6).Value = "2.22"
$ws.Cells.Item(14,7).Value = "0.0056"
$ws.Cells.Item(14,8).Value = 7

$ws.Cells.Item(15,1).Value = 13
$ws.Cells.Item(15,2).Value = "003855"
$ws.Cells.Item(15,3).Value = "汇安丰华灵活配置混合C"
$ws.Cells.Item(15,4).Value = "0.19"
$ws.Cells.Item(15,5).Value = "45.55"
$ws.Cells.Item(15,6).Value = "2.33"
$ws.Cells.Item(15,7).Value = "0.0044"
$ws.Cells.Item(15,8).Value = 6

$ws.Cells.Item(16,1).Value = 14
$ws.Cells.Item(16,2).Value = "003854"
$ws.Cells.Item(16,3).Value = "汇安丰华灵活配置混合A"
$ws.Cells.Item(16,4).Value = "0.00"
$ws.Cells.Item(16,5).Value = "45.55"
$ws.Cells.Item(16,6).Value = "2.33"
$ws.Cells.Item(16,7).Value = 0
$ws.Cells.Item(16,8).Value = 6

$indexStyleSrc.Copy()
$ws.Range("A2:A16").PasteSpecial(-4122)

# Force the numeric-looking D/E/F (and G2:G15) columns to stay text, as in
# the source data (keeps literal formatting such as trailing zeros, e.g.
# "35.14" / "0.00"), then clear the number-format style back off the cells
# so they fall back onto the default style (matches the rest of the sheet,
# which has no explicit style on the data columns).
$textRange = $ws.Range("D2:G16")
$textRange.NumberFormat = "@"
$ws.Cells.Item(2,4).Value = "35.14"
$ws.Cells.Item(2,5).Value = "90.17"
$ws.Cells.Item(2,6).Value = "6.69"
$ws.Cells.Item(2,7).Value = "2.3509"
$ws.Cells.Item(3,4).Value = "4.46"
$ws.Cells.Item(3,5).Value = "90.97"
$ws.Cells.Item(3,6).Value = "9.04"
$ws.Cells.Item(3,7).Value = "0.4032"
$ws.Cells.Item(4,4).Value = "3.85"
$ws.Cells.Item(4,5).Value = "98.34"
$ws.Cells.Item(4,6).Value = "3.87"
$ws.Cells.Item(4,7).Value = "0.1490"
$ws.Cells.Item(5,4).Value = "2.57"
$ws.Cells.Item(5,5).Value = "81.91"
$ws.Cells.Item(5,6).Value = "3.06"
$ws.Cells.Item(5,7).Value = "0.0786"
$ws.Cells.Item(6,4).Value = "2.65"
$ws.Cells.Item(6,5).Value = "74.01"
$ws.Cells.Item(6,6).Value = "2.22"
$ws.Cells.Item(6,7).Value = "0.0588"
$ws.Cells.Item(7,4).Value = "1.47"
$ws.Cells.Item(7,5).Value = "85.46"
$ws.Cells.Item(7,6).Value = "3.96"
$ws.Cells.Item(7,7).Value = "0.0582"
$ws.Cells.Item(8,4).Value = "0.96"
$ws.Cells.Item(8,5).Value = "52.43"
$ws.Cells.Item(8,6).Value = "5.52"
$ws.Cells.Item(8,7).Value = "0.0530"
$ws.Cells.Item(9,4).Value = "0.87"
$ws.Cells.Item(9,5).Value = "92.08"
$ws.Cells.Item(9,6).Value = "3.88"
$ws.Cells.Item(9,7).Value = "0.0338"
$ws.Cells.Item(10,4).Value = "0.34"
$ws.Cells.Item(10,5).Value = "92.91"
$ws.Cells.Item(10,6).Value = "7.81"
$ws.Cells.Item(10,7).Value = "0.0266"
$ws.Cells.Item(11,4).Value = "0.44"
$ws.Cells.Item(11,5).Value = "52.43"
$ws.Cells.Item(11,6).Value = "5.52"
$ws.Cells.Item(11,7).Value = "0.0243"
$ws.Cells.Item(12,4).Value = "0.27"
$ws.Cells.Item(12,5).Value = "52.43"
$ws.Cells.Item(12,6).Value = "5.52"
$ws.Cells.Item(12,7).Value = "0.0149"
$ws.Cells.Item(13,4).Value = "0.36"
$ws.Cells.Item(13,5).Value = "85.46"
$ws.Cells.Item(13,6).Value = "3.96"
$ws.Cells.Item(13,7).Value = "0.0143"
$ws.Cells.Item(14,4).Value = "0.25"
$ws.Cells.Item(14,5).Value = "74.01"
$ws.Cells.Item(14,6).Value = "2.22"
$ws.Cells.Item(14,7).Value = "0.0056"
$ws.Cells.Item(15,4).Value = "0.19"
$ws.Cells.Item(15,5).Value = "45.55"
$ws.Cells.Item(15,6).Value = "2.33"
$ws.Cells.Item(15,7).Value = "0.0044"
$ws.Cells.Item(16,4).Value = "0.00"
$ws.Cells.Item(16,5).Value = "45.55"
$ws.Cells.Item(16,6).Value = "2.33"
# G16 is a genuine number (0) in the source data, unlike G2:G15, so restore
# it after the text coercion above.
$ws.Cells.Item(16,7).Value = 0
$textRange.Style = "Normal"

$ws.Range("A1").Select()

# ---------------------------------------------------------------------------
# 2. Record the new quarter in the "总计" (totals) summary sheet: insert a
#    new row 2 (shifting the existing quarters down) with the 2022-Q3 totals.
# ---------------------------------------------------------------------------
$totalIndexStyleSrc = $totalSheet.Range("A2")

$totalSheet.Cells.Item(9,1).Value = 7
$totalSheet.Cells.Item(9,2).Value = "2020-Q4"
$totalSheet.Cells.Item(9,3).Value = 9
$totalSheet.Cells.Item(9,4).Value = 1.99

$totalSheet.Cells.Item(8,1).Value = 6
$totalSheet.Cells.Item(8,2).Value = "2021-Q1"
$totalSheet.Cells.Item(8,3).Value = 4
$totalSheet.Cells.Item(8,4).Value = 3.36

$totalSheet.Cells.Item(7,1).Value = 5
$totalSheet.Cells.Item(7,2).Value = "2021-Q2"
$totalSheet.Cells.Item(7,3).Value = 3
$totalSheet.Cells.Item(7,4).Value = 0.73

$totalSheet.Cells.Item(6,1).Value = 4
$totalSheet.Cells.Item(6,2).Value = "2021-Q3"
$totalSheet.Cells.Item(6,3).Value = 4
$totalSheet.Cells.Item(6,4).Value = 0.68

$totalSheet.Cells.Item(5,1).Value = 3
$totalSheet.Cells.Item(5,2).Value = "2021-Q4"
$totalSheet.Cells.Item(5,3).Value = 7
$totalSheet.Cells.Item(5,4).Value = 3.46

$totalSheet.Cells.Item(4,1).Value = 2
$totalSheet.Cells.Item(4,2).Value = "2022-Q1"
$totalSheet.Cells.Item(4,3).Value = 16
$totalSheet.Cells.Item(4,4).Value = 8.72

$totalSheet.Cells.Item(3,1).Value = 1
$totalSheet.Cells.Item(3,2).Value = "2022-Q2"
$totalSheet.Cells.Item(3,3).Value = 9
$totalSheet.Cells.Item(3,4).Value = 3.19

$totalSheet.Cells.Item(2,1).Value = 0
$totalSheet.Cells.Item(2,2).Value = "2022-Q3"
$totalSheet.Cells.Item(2,3).Value = 15
$totalSheet.Cells.Item(2,4).Value = 3.28

$totalIndexStyleSrc.Copy()
$totalSheet.Range("A2:A9").PasteSpecial(-4122)

$totalSheet.Range("A1").Select()
$totalSheet.Activate()
